# Auto-update draw results: append the 2025-10-06 Pick 3 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All columns in this sheet store plain text (dates, zero-padded codes, and
# the dash-separated result are never real numbers), so force text formatting
# on the new row before writing values - otherwise Excel would silently
# reinterpret strings like "2025-10-06" or "251006" as a date/number.
$ws.Range("A20:E20").NumberFormat = "@"

$ws.Range("A20").Value = "2025-10-06"
$ws.Range("B20").Value = "Pick 3"
$ws.Range("C20").Value = "251006"
$ws.Range("D20").Value = "7-8-8"
$ws.Range("E20").Value = "2025-10-06T21:37:29.248+04:00"
